# Regenerate the "K" (strikeout) column (G) with newly calculated values.
# This replaces the previous placeholder "Strike#" based values that had
# been written into column G with the real K values computed from the
# refreshed save_data pull.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 1
    4  = 2
    5  = 2
    6  = 0
    8  = 0
    9  = 2
    10 = 3
    11 = 1
    12 = 3
    13 = 1
    14 = 1
    15 = 1
    16 = 3
    17 = 0
    18 = 1
    19 = 1
    20 = 0
    21 = 3
    22 = 0
    23 = 2
    24 = 4
    26 = 3
    27 = 1
    29 = 2
    30 = 2
    31 = 0
    32 = 1
    33 = 1
    34 = 2
    35 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
